# Actualizacion Datos Personales 4 nov
#
# Student "VAZQUEZ VERA MARIA JOSE" (row 45 on "Calificaciones") gets a real
# grade of 5 for "TECNOLOGIAS DE LA INFORMACION Y LA COMUNICACION" in both the
# 1P column (F) and the Final column (X) - it used to be blank (-1).
#
# That removes one "blank" record for her on the "Blancos" sheet (row 65,
# the TECNOLOGIAS entry), which shifts the remaining rows up by one and
# shrinks the used range by a row. It also changes her blanks-count on
# "Totales Blanco" from 4 to 3, which re-sorts that block of rows (now tied
# with two other students at 3 blanks each, ordered by ascending NC). And it
# updates the aggregated stats for that subject/teacher on "Totales" (row 7):
# one more reprobado, one less blanco.

$wb = $excel.ActiveWorkbook

# --- 1) Calificaciones: fill in the previously-blank grade -----------------
$wsCal = $wb.Worksheets.Item("Calificaciones")
$wsCal.Range("F45").Value = 5
$wsCal.Range("X45").Value = 5

# --- 2) Totales: row 7 (TECNOLOGIAS / Pesce Bautista Victor Manuel) --------
$wsTot = $wb.Worksheets.Item("Totales")
$wsTot.Range("E7").Value = 1          # Reprobados: 0 -> 1
$wsTot.Range("G7").Value = 2.33       # Por_Repro:  0 -> 2.33
$wsTot.Range("H7").Value = 8.5        # Promedio:   8.6 -> 8.5
$wsTot.Range("I7").Value = 4          # Blancos:    5 -> 4
$wsTot.Range("J7").Value = 9.300000000000001   # Por_Blancos: 11.63 -> 9.3...

# --- 3) Blancos: remove the now-resolved blank record (row 65) -------------
# Row 65 was: 21330051920109 / VAZQUEZ / VERA / MARIA JOSE / TECNOLOGIAS.. / Pesce
# Deleting it shifts rows 66-67 up to become rows 65-66, and shrinks the
# worksheet's used range from A1:F67 down to A1:F66.
$wsBlancos = $wb.Worksheets.Item("Blancos")
$wsBlancos.Rows.Item(65).Delete()

# --- 4) Totales Blanco: re-sort rows 9-11 after the blanks count change ----
# Row 9 (VAZQUEZ, 21330051920109) drops from 4 blanks to 3, tying it with
# the two rows below; the three tied rows end up ordered by ascending NC.
$wsTotBlanco = $wb.Worksheets.Item("Totales Blanco")
$wsTotBlanco.Range("A9").Value = 21330051920089
$wsTotBlanco.Range("B9").Value = "MORALES"
$wsTotBlanco.Range("C9").Value = "RODRIGUEZ"
$wsTotBlanco.Range("D9").Value = "FATIMA MARILYN"
$wsTotBlanco.Range("E9").Value = 3

$wsTotBlanco.Range("A10").Value = 21330051920096
$wsTotBlanco.Range("B10").Value = "PEREZ"
$wsTotBlanco.Range("C10").Value = "PAZ"
$wsTotBlanco.Range("D10").Value = "JAIME"
$wsTotBlanco.Range("E10").Value = 3

$wsTotBlanco.Range("A11").Value = 21330051920109
$wsTotBlanco.Range("B11").Value = "VAZQUEZ"
$wsTotBlanco.Range("C11").Value = "VERA"
$wsTotBlanco.Range("D11").Value = "MARIA JOSE"
$wsTotBlanco.Range("E11").Value = 3
